# Update benchmark: 2026-01-31 06:47:59 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (HESAPTAN EFT - Şube): clear F3
$ws.Range("F3").Value = $null

# Row 4 (HESAPTAN EFT - ATM): clear F4
$ws.Range("F4").Value = $null

# Row 5 (HESAPTAN EFT - Mobil): clear F5
$ws.Range("F5").Value = $null

# Row 6 (DÜZENLİ EFT): set H6
$ws.Range("H6").Value = "8.300,01 TL - 199,41 TL"

# Row 8 (HESAPTAN HAVALE - Şube): clear F8
$ws.Range("F8").Value = $null

# Row 9 (HESAPTAN HAVALE - ATM): clear F9
$ws.Range("F9").Value = $null

# Row 10 (HESAPTAN HAVALE - Mobil): clear F10
$ws.Range("F10").Value = $null

# Row 13 (GELEN SWIFT)
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"
$ws.Range("F13").Value = $null
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 7,97 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14 (GİDEN SWIFT - Mobil)
$ws.Range("F14").Value = $null
$ws.Range("H14").Value = "3.000 TL - 6.000 TL"
